$wb = $excel.ActiveWorkbook
$wsMoorings = $wb.Worksheets.Item("Moorings")
$wsAsset    = $wb.Worksheets.Item("Asset_Cal_Info")

# ---------------------------------------------------------------------------
# Moorings sheet ("sheet1") - glider GL001 -> GL276, deployment 2 -> 1
# ---------------------------------------------------------------------------
$wsMoorings.Range("A2").Value = "GP05MOAS-GL276"
$wsMoorings.Range("C2").Value = 1

$wsMoorings.Range("G2").Value = "50° 04.70' N"
$wsMoorings.Range("H2").Value = "144° 48.32' W"
$wsMoorings.Range("I2").Value = "0m"
$wsMoorings.Range("J2").Value = "MV1404"

$wsMoorings.Range("L1").Value = "Lat"
$wsMoorings.Range("M1").Value = "Lon"
$wsMoorings.Range("N1").Value = "Data Start"
$wsMoorings.Range("O1").Value = "Data End"

# ---------------------------------------------------------------------------
# Asset_Cal_Info sheet ("sheet2") - reference designators + deployment number
# ---------------------------------------------------------------------------
$wsAsset.Range("A3").Value = "GP05MOAS-GL276-00-ENG000000"
$wsAsset.Range("C3").Value = 1
$wsAsset.Range("G3").Value = "OpenOceanGlider_276_Factory_Configs_Calibrations_2013-10-01.pdf"

$wsAsset.Range("A4").Value = "GP05MOAS-GL276-01-FLORDM000"
$wsAsset.Range("C4").Value = 1
$wsAsset.Range("G4").Value = "Default value per <flo_bback_total(beta, degC=20.0, psu=32.0, theta=117.0, wlngth=700.0, xfactor=1.08)>"

$wsAsset.Range("A5").Value = "GP05MOAS-GL276-01-FLORDM000"
$wsAsset.Range("C5").Value = 1
$wsAsset.Range("G5").Value = "Default value per <flo_bback_total(beta, degC=20.0, psu=32.0, theta=117.0, wlngth=700.0, xfactor=1.08)>"

$wsAsset.Range("A6").Value = "GP05MOAS-GL276-01-FLORDM000"
$wsAsset.Range("C6").Value = 1
$wsAsset.Range("G6").Value = "Default value per <flo_bback_total(beta, degC=20.0, psu=32.0, theta=117.0, wlngth=700.0, xfactor=1.08)>"

$wsAsset.Range("A7").Value = "GP05MOAS-GL276-01-FLORDM000"
$wsAsset.Range("C7").Value = 1
$wsAsset.Range("G7").Value = "Default value per <flo_scat_seawater(degC, psu, theta=117.0, wlngth=700.0, delta=0.039)>"

$wsAsset.Range("A8").Value = "GP05MOAS-GL276-02-DOSTAM000"
$wsAsset.Range("C8").Value = 1
$wsAsset.Range("G8").Value = "Requires Lat, Lon, pressure, and temperature from glider engineering (PD1382/1391 1527/1528) and PRACSAL_L2 (PD1560); AADI DOSTA - Anderaa Optode"

$wsAsset.Range("A9").Value = "GP05MOAS-GL276-04-CTDGVM000"
$wsAsset.Range("C9").Value = 1
$wsAsset.Range("G9").Value = "Requires Lat and Lon from glider engineering (PD1382/1391); Seabird Pumped CTD"

# ---------------------------------------------------------------------------
# View state: active sheet moves from Moorings to Asset_Cal_Info, and the
# selections on each sheet move too.
# ---------------------------------------------------------------------------
$wsMoorings.Activate() | Out-Null
$wsMoorings.Range("A2").Select() | Out-Null

$wsAsset.Activate() | Out-Null
$wsAsset.Range("C10").Select() | Out-Null
